$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("my worksheet")

# Names in column A
$names = @("John", "Denise", "Adam", "Jasmine", "Phoebe ", "Gryff", "Milton", "Gingy", "Joppa")
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 1).Value = $names[$i]
}

# Column B values
$colB = @(21, 18, 23, 20, 22, 23, 17, 16, 36)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 2).Value = $colB[$i]
}

# Column C values
$colC = @(56, 42, 21, 20, 22, 23, 67, 54, 35)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 3).Value = $colC[$i]
}

# Column D values
$colD = @(21, 34, 32, 45, 55, 67, 34, 87, 56)
for ($i = 0; $i -lt $colD.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 4).Value = $colD[$i]
}

# Column E: average formula, row 1 standalone, rows 2-9 as one filled/shared range
$ws1.Range("E1").Formula = "=(C1+D1)/2"
$ws1.Range("E2:E9").Formula = "=(C2+D2)/2"

$ws1.Range("F1").Select() | Out-Null
